# CARAGA_NEWCON.xlsx - refresh to the most recently updated status /
# accomplishment report (as of May): retire the old "No. of ... (Sites/CL
# status breakdown)" helper columns (AB:AS) in favor of a single
# "Program-Year" column, rename the classroom-count header, and update the
# first data row's category / program-year values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "No. of Classrooms" -> two-line "PHYSICAL TARGET / (# OF CL)" header
$ws.Range("I1").Value = "PHYSICAL TARGET" + [char]10 + "(# OF CL)"

# "No. of Sites Reverted" -> "Program-Year" (new header for column AB)
$ws.Range("AB1").Value = "Program-Year"

# Drop the old per-status breakdown columns (AC:AS) entirely. This shifts
# the trailing "Status as of July 4, 2025" column (previously AT) left into
# AC, and the dependent data validation's sqref moves with it.
$ws.Range("AC1:AS2").EntireColumn.Delete()

# Row 2 (Falcon Memorial ES) updates to the latest accomplishment data
$ws.Range("A2").Value = "K10"
$ws.Range("AB2").Value = "NC 2020"
